# edit.ps1 - applies Lab05 report content changes via Word COM-interop
$d = $word.ActiveDocument

function Get-XmlFromBase64($b64) {
    $bytes = [System.Convert]::FromBase64String($b64)
    return [System.Text.Encoding]::UTF8.GetString($bytes)
}

$titleXml   = Get-XmlFromBase64 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6dD5MYWIgMDUgUmVwb3J0IOKAkyBVTjwvdzp0PjwvdzpyPjx3OnI+PHc6dD5ldCBPdmVydmlldzwvdzp0PjwvdzpyPjwvdzpwPg=="
$introXml   = Get-XmlFromBase64 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5XZSB0b29rIHRoZSB0aW1lIG9uIHRoaXMgbGFiIHRvIHdhdGNoIGFuZCBmb2xsb3cgdGhlIGluc3RydWN0aW9ucyBvZiBZb3V0dWJlciA8L3c6dD48L3c6cj48dzpyPjx3OnJQcj48dzppLz48L3c6clByPjx3OnQ+R2FtZXIgdG8gR2FtZSBEZXZlbG9wZXI8L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IG92ZXIgYSB0d28gcGFydCB2aWRlbyB3aGVyZSBoZSBleHBsYWlucyB0aGUgYmFzaWNzIGJlaGluZCB0aGUgbmV3IG5ldHdvcmtpbmcgc3VpdGUgVU5ldC4gSGUgZXhwbGFpbnMgdGhlIHN0ZXBzIHRvIHF1aWNrbHkgY29ubmVjdCBzZXNzaW9ucyBvdmVyIGEgbG9jYWxob3N0LiBIZSB0aGVuIHJlbW92ZXMgc29tZSBjb21wb25lbnRzIHNvIGhlIGNhbiBzaG93IHVzIGEgc21vb3RoZXIsIGVhc2llciwgYW5kIG1vcmUgcm9idXN0IHZlcnNpb24gb2Ygd2hhdCB3YXMgYWxyZWFkeSBhbGxvd2VkLjwvdzp0PjwvdzpyPjwvdzpwPg=="
$methodsXml = Get-XmlFromBase64 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6dD5IZSBzdGFydHMgb2ZmIGJ5IGhhdmluZyB1cyB1c2UgdGhlIHN0YW5kYXJkIGFzc2V0cyBwYWNrYWdlIGZyb20gdGhlIHVuaXR5IHN0b3JlIGFuZCB0aGUgYmFzZSBlbmdpbmUgdG8gY3JlYXRlIGEgZmV3IGNvcmUgY29tcG9uZW50cyB0byBtYW5pcHVsYXRlOiBhIHBsYW5lIGZvciB0cmF2ZXJzYWwsIHR3byBjdWJlcyBvZiBkaWZmZXJpPC93OnQ+PC93OnI+PHc6cj48dzp0Pm5nIHNpemUsIGFuZCBhIGZpcnN0IHBlcnNvbiBjb250cm9sbGVyIG5hbWVkIHBsYXllciB3aXRoIGFuIGFkZGVkIGNhcHN1bGUgZm9yIHZpc3VhbC4gSGUgYWRkZWQgYSBxdWljayBjaGFuZ2UgdG8gdGhlIHBsYW5lIGJ5IGFkZGluZyBhIHNpbXBsZSBncmF5IG1hdGVyaWFsIGFuZCB0aGVuIHNldCBhbGwgYnV0IHRoZSBwbGF5ZXIgaW4gYW4gZW1wdHkgZ2FtZSBvYmplY3QgdG8gcmVtb3ZlIGNsdXR0ZXIgZnJvbSB0aGUgSGllcmFyY2h5Ljwvdzp0PjwvdzpyPjwvdzpwPjx3OnAgeG1sbnM6dz0iaHR0cDovL3NjaGVtYXMub3BlbnhtbGZvcm1hdHMub3JnL3dvcmRwcm9jZXNzaW5nbWwvMjAwNi9tYWluIj48dzpyPjx3OnQ+Tm93IHdlIGdldCB0byB0aGUgbWVhdCBvZiB0aGUgaXNzdWUuICBIZSBjcmVhdGVkIGEgYmxhbmsgb2JqZWN0IGNhbGxlZCBOZXR3b3JrTWFuYWdlciBhbmQgc2VhcmNoZWQgZm9yIGEgTmV0d29yayBNYW5hZ2VyIGNvbXBvbmVudCBmcm9tIHRoZSBBZGQgY29tcG9uZW50IHRhYi4gIEhlIHRoZW4gdXNlZCB0aGUgc2FtZSB0ZWNobmlxdWUgdG8gYWRkIGEgTmV0d29yayBNYW5hZ2VyIEhVRDwvdzp0PjwvdzpyPjx3OnI+PHc6dD4uPC93OnQ+PC93OnI+PC93OnA+PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj5XZSBtYWtlIGEgcHJlZmFiIG9mIHRoZSBwbGF5ZXIgYW5kIHRoZW4gYWRkIHRoYXQgcHJlZmFiIHRvIHRoZSDigJxQbGF5ZXIgUHJlZmFi4oCdIG9wdGlvbiBvZiB0aGUgTmV0d29yayBNYW5hZ2VyIGNvbXBvbmVudCBzbyBpdCBpcyBhZGRlZCB3aGVuIHRoZSBwZXJzb24gcGxheWluZyBpcyBhZGRlZCB0byB0aGUgc2NlbmUuIEhlIHRoZW4gYWRkcyB0d28gbmV3IGNvbXBvbmVudHMgdG8gdGhlIHBsYXllcjogdGhlIE5ldHdvcmsgSWRlbnRpdHkgYW5kIE5ldHdvcmsgVHJhbnNmb3JtIHNjcmlwdHMuIEhlIGFjdGl2YXRlcyB0aGUgPC93OnQ+PC93OnI+PHc6cj48dzp0PkxvY2FsIFBsYXllciBBdXRob3JpdHkgbm9kZSBvZiB0aGUgTmV0d29yayBUcmFuc2Zvcm0gdG8gaGF2ZSB0aGUgYWN0aXZlIHBsYXllciBvZiB0aGUgcGMgdXNlIGl0IGFuZCBzZXRzIHRoZSBUcmFuc2Zvcm0gU3luYyBNb2RlIHRvPC93OnQ+PC93OnI+PHc6cj48dzp0IHhtbDpzcGFjZT0icHJlc2VydmUiPiDigJxTeW5jIFRyYW5zZm9ybeKAnTwvdzp0PjwvdzpyPjx3OnI+PHc6dCB4bWw6c3BhY2U9InByZXNlcnZlIj4gb24gdGhlIE5ldHdvcmsgVHJhbnNmb3JtIGluc3RlYWQgb2YgdGhlIGRlZmF1bHQgc28gdGhlIHBsYXllcnMgdHJhbnNmb3JtIGlzIGNvbW11bmljYXRlZC48L3c6dD48L3c6cj48dzpib29rbWFya1N0YXJ0IHc6aWQ9IjAiIHc6bmFtZT0iX0dvQmFjayIvPjx3OmJvb2ttYXJrRW5kIHc6aWQ9IjAiLz48L3c6cD48dzpwIHhtbG5zOnc9Imh0dHA6Ly9zY2hlbWFzLm9wZW54bWxmb3JtYXRzLm9yZy93b3JkcHJvY2Vzc2luZ21sLzIwMDYvbWFpbiI+PHc6cj48dzp0PkFmdGVyIHdlIG1ha2Ugb3VyIGZpcnN0IGJ1aWxk4oCUYWRkaW5nIGJvdGggdGhlIG9mZmxpbmUgbWVudSBzY2VuZSBhbmQgb25saW5lIG1haW4gc2NlbmXigJR3ZSBkcmFnIHRoZSBtZW51IGFuZCBtYWluIHNjZW5lcyBpbnRvIHRoZSBvZmZsaW5lIGFuZCBvbmxpbmUgc2NlbmVzIHJlc3BlY3RpdmVseS4gVGhpcyBzd2l0Y2hlcyBiZXR3ZWVuIHRoZSB0d28gc2NlbmVzIGRlcGVuZGluZyBvbiB0aGUgb25saW5lIHN0YXR1cy48L3c6dD48L3c6cj48dzpyPjx3OnQgeG1sOnNwYWNlPSJwcmVzZXJ2ZSI+IDwvdzp0PjwvdzpyPjwvdzpwPg=="
$emptyXml   = Get-XmlFromBase64 "PHc6cCB4bWxuczp3PSJodHRwOi8vc2NoZW1hcy5vcGVueG1sZm9ybWF0cy5vcmcvd29yZHByb2Nlc3NpbmdtbC8yMDA2L21haW4iLz4="

# --- 1. Title paragraph: split "Unet" (spell-checked) into "UN" + "et" runs ---
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Lab 05 Report*Overview*") {
        $titlePara = $p
        break
    }
}
$r = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$null = $r.InsertXML($titleXml)

# --- 2. Introduction body paragraph: empty paragraph right after "Introduction" ---
$introHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Introduction`r") {
        $introHeading = $p
        break
    }
}
$introBody = $introHeading.Next()
$r = $d.Range($introBody.Range.Start, $introBody.Range.Start)
$null = $r.InsertXML($introXml)

# --- 3. Methods and Work: empty paragraph right after heading expands to 4 paragraphs ---
$methodsHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Methods and Work`r") {
        $methodsHeading = $p
        break
    }
}
$methodsBody = $methodsHeading.Next()
$r = $d.Range($methodsBody.Range.Start, $methodsBody.Range.Start)
$null = $r.InsertXML($methodsXml)

# --- 4. Remove old trailing paragraph with the "_GoBack" bookmark (now just an empty paragraph) ---
$lastPara = $d.Paragraphs.Last
$r = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$null = $r.InsertXML($emptyXml)

Write-Output "Edit complete."
